$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure new comments are attributed to "Narongdej Sarnsuwan"
$excel.UserName = "Narongdej Sarnsuwan"

# --- Cell values ---
# Shared-string table order matters: "wealth" must be inserted before
# "awesome" so they land at uniqueCount indexes 1 and 2 respectively.
$ws.Range("A3").Value = "wealth"
$ws.Range("A1").Value = "awesome"
$ws.Range("B3").Value = 200

# --- Comment on B3 ---
$commentText = "Narongdej Sarnsuwan:" + "`n" + "{{--" + "`n" + "key: <<A1>>_wealth" + "`n" + "fallback: <<A5>>" + "`n" + "--}}"
$ws.Range("B3").AddComment($commentText)

# --- Selection ---
$ws.Range("L15").Select()
